# S19 Table update: revised metabolite counts / recomputed percentages,
# updated title text (208 -> 206 metabolites), refreshed number format
# for the Percentage column, and one extra (blank) trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Title cell (A1): replace "208" with "206", keep the existing
#        italic run on "Lactobacillus plantarum" intact. ---
$titleCell = $ws.Range("A1")
$titleCell.Characters(50, 3).Text = "206"
$fullTitle = $titleCell.Value2
$speciesIdx = $fullTitle.IndexOf("Lactobacillus")
$speciesLen = ("Lactobacillus plantarum").Length
$speciesRun = $titleCell.Characters($speciesIdx + 1, $speciesLen)
$speciesRun.Font.Italic = $true
$speciesRun.Font.Name = "Calibri"
$speciesRun.Font.Size = 11

# --- 2. Updated frequency (column B) and recomputed percentage
#        (column C) values, row by row. ---
$ws.Range("B4").Value = 153
$ws.Range("C4").Value = 23.042168674698797

$ws.Range("B5").Value = 146
$ws.Range("C5").Value = 21.987951807228917

$ws.Range("B6").Value = 62
$ws.Range("C6").Value = 9.3373493975903621

$ws.Range("B7").Value = 51
$ws.Range("C7").Value = 7.6807228915662646

$ws.Range("B8").Value = 47
$ws.Range("C8").Value = 7.0783132530120483

$ws.Range("B9").Value = 46
$ws.Range("C9").Value = 6.927710843373494

$ws.Range("B10").Value = 27
$ws.Range("C10").Value = 4.0662650602409638

$ws.Range("B11").Value = 24
$ws.Range("C11").Value = 3.6144578313253013

$ws.Range("B12").Value = 21
$ws.Range("C12").Value = 3.1626506024096384

$ws.Range("B13").Value = 16
$ws.Range("C13").Value = 2.4096385542168677

$ws.Range("B14").Value = 13
$ws.Range("C14").Value = 1.9578313253012047

$ws.Range("B15").Value = 12
$ws.Range("C15").Value = 1.8072289156626506

$ws.Range("B16").Value = 11
$ws.Range("C16").Value = 1.6566265060240963

$ws.Range("B17").Value = 11
$ws.Range("C17").Value = 1.6566265060240963

$ws.Range("B18").Value = 7
$ws.Range("C18").Value = 1.0542168674698795

$ws.Range("B19").Value = 6
$ws.Range("C19").Value = 0.90361445783132532

$ws.Range("B20").Value = 6
$ws.Range("C20").Value = 0.90361445783132532

$ws.Range("B21").Value = 3
$ws.Range("C21").Value = 0.45180722891566266

$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 0.15060240963855423

$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 0.15060240963855423

# --- 3. Columns A & B (rows 4-23) drop their old centered style,
#        reverting to the default/general formatting. ---
$ws.Range("A4:B23").ClearFormats()

# --- 4. Column C (rows 4-24) gets the refreshed "0.0" number format
#        (replacing the old centered "0.00" style) with no special
#        alignment applied. ---
$ws.Range("C4:C24").ClearFormats()
$ws.Range("C4:C24").NumberFormat = "0.0"

# --- 5. New trailing blank row 24 (only C24 carries the number
#        format; no value) is created implicitly because C24 is part
#        of the C4:C24 NumberFormat range above. ---

# --- 6. Selection cosmetic update, matching the saved workbook state. ---
$ws.Range("B4").Select()
